# Refresh the cryptos price/volume snapshot (scheduled GitHub Actions update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds plain text in the source sheet (t="inlineStr"), even when the
# text happens to look like a plain number (e.g. "562.21"). Assigning such a string via
# .Value would otherwise be auto-coerced to a numeric cell, so those cells are switched to
# the Text number format first to preserve them as text, matching the original file.

$ws.Range("D2").Value = "68.915.71"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "2.475.33"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.21"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.70"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.158"
$ws.Range("E9").Value = "  +4.98%  "
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.86"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").Value = "68.845.30"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.70"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "10.68"
$ws.Range("E16").Value = "  -2.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "339.02"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.97"
$ws.Range("E18").Value = "  -2.91%  "
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.90"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.93"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.69"
$ws.Range("E23").Value = "  -0.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.29"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").Value = "0.0₃0827"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.26"
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "432.80"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.96"
$ws.Range("E31").Value = "  +1.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.02"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.109"
$ws.Range("E34").Value = "  -1.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.92"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.48"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.300"
$ws.Range("E37").Value = "  -2.02%  "
$ws.Range("E38").Value = "  -2.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.09"
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.08"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.40"
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "131.00"
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0721"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.488"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0921"
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.03"
$ws.Range("E49").Value = "  -6.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.95"
$ws.Range("E50").Value = "  -3.55%  "
$ws.Range("D51").Value = "0.0₆0208"
$ws.Range("E51").Value = "  -7.52%  "
